$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.370.28'
$ws.Range("E2").Value = '  -2.74%  '
$ws.Range("D3").Value = '1.859.01'
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '330.68'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").Value = '0.4736'
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("D8").Value = '0.3963'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("D9").Value = '47.16'
$ws.Range("E9").Value = '  -11.06%  '
$ws.Range("D10").Value = '0.08003'
$ws.Range("E10").Value = '  -4.85%  '
$ws.Range("D11").Value = '1.018'
$ws.Range("E11").Value = '  -2.55%  '
$ws.Range("D12").Value = '21.55'
$ws.Range("E12").Value = '  -2.62%  '
$ws.Range("D13").Value = '1.853.77'
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.960'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("D15").Value = '7.166'
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '86.41'
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("D18").Value = '0.00001037'
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("D19").Value = '0.06545'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '17.23'
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '5.496'
$ws.Range("E22").Value = '  -4.19%  '
$ws.Range("D23").Value = '27.369.40'
$ws.Range("E23").Value = '  -2.77%  '
$ws.Range("D24").Value = '10.94'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("D25").Value = '2.298'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = '2.074.39'
$ws.Range("E26").Value = '  -3.46%  '
$ws.Range("D27").Value = '20.39'
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("D28").Value = '154.12'
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("D29").Value = '2.077'
$ws.Range("E29").Value = '  -2.83%  '
$ws.Range("D30").Value = '5.507'
$ws.Range("E30").Value = '  -4.44%  '
$ws.Range("D31").Value = '122.29'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").Value = '0.09507'
$ws.Range("E32").Value = '  -1.39%  '
$ws.Range("D33").Value = '0.9553'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = '1.447'
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("D35").Value = '3.589'
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("D36").Value = '5.285'
$ws.Range("E36").Value = '  -4.66%  '
$ws.Range("D37").Value = '0.06059'
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("D38").Value = '0.02236'
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("D39").Value = '1.209'
$ws.Range("E39").Value = '  -4.38%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = '8.022'
$ws.Range("E41").Value = '  -9.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5950'
$ws.Range("E42").Value = '  -3.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1900'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '10.29'
$ws.Range("E44").Value = '  -6.63%  '
$ws.Range("D45").Value = '1.267'
$ws.Range("E45").Value = '  -3.52%  '
$ws.Range("D46").Value = '0.5663'
$ws.Range("E46").Value = '  -3.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.10'
$ws.Range("E47").Value = '  -5.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.430'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.930'
$ws.Range("E49").Value = '  -4.57%  '
$ws.Range("D50").Value = '0.06767'
$ws.Range("E50").Value = '  -2.06%  '
$ws.Range("D51").Value = '109.51'
$ws.Range("E51").Value = '  -1.64%  '
